$d = $word.ActiveDocument

# 1. Location in the dateline: Valparaíso -> Santiago
$d.Content.Find.Execute("Valparaíso", $true, $false, $false, $false, $false, $true, 1, $false, "Santiago", 2) | Out-Null

# 2. Day of the month in the dateline: 24 -> 28
$d.Content.Find.Execute("24", $true, $false, $false, $false, $false, $true, 1, $false, "28", 2) | Out-Null

# 3. Signatory name: Eduardo Muñoz Inchausti -> Jeannette Rodríguez Chandia
$d.Content.Find.Execute("Eduardo Muñoz Inchausti", $true, $false, $false, $false, $false, $true, 1, $false, "Jeannette Rodríguez Chandia", 2) | Out-Null

# 4. Signatory title: Director -> Jefa de Carrera Campus Santiago
$d.Content.Find.Execute("Director", $true, $false, $false, $false, $false, $true, 1, $false, "Jefa de Carrera Campus Santiago", 2) | Out-Null

# 5. Initials: EMI -> JRC
$d.Content.Find.Execute("EMI", $true, $false, $false, $false, $false, $true, 1, $false, "JRC", 2) | Out-Null

# 6. Typist initials: krr -> lll
$d.Content.Find.Execute("krr", $true, $false, $false, $false, $false, $true, 1, $false, "lll", 2) | Out-Null

# 7. Footer contact line -> new Campus Santiago address/phone
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("Las Heras Nº 06 Valparaíso | Fono: (32) 250 7961- 2507815 | E-mail: practivasv@uv.cl, www.uv.cl", $true, $false, $false, $false, $false, $true, 1, $false, "Campus Santiago - Gran Avenida 4160, San Miguel | Fono +56 (2)2329  2149", 2) | Out-Null
